$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper pattern: capture style, force text format so numeric-looking
# strings (e.g. "1.000", "242.18") are preserved verbatim as text instead
# of being auto-coerced to a number, then restore the original style so
# no stray style/format is left behind on the cell.

$ws.Range("D2").Value = "30.030.44"
$ws.Range("E2").Value = "  -0.92%  "
$ws.Range("D3").Value = "1.902.13"
$ws.Range("E3").Value = "  -1.44%  "
$c = $ws.Range("D4")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "1.000"
$c.Style = $s
$ws.Range("E4").Value = "  -0.04%  "
$c = $ws.Range("D5")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.7442"
$c.Style = $s
$ws.Range("E5").Value = "  -0.61%  "
$c = $ws.Range("D6")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "242.18"
$c.Style = $s
$ws.Range("E6").Value = "  -0.64%  "
$c = $ws.Range("D7")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "1.000"
$c.Style = $s
$ws.Range("E7").Value = "  -0.02%  "
$c = $ws.Range("D8")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.3070"
$c.Style = $s
$ws.Range("E8").Value = "  -3.34%  "
$c = $ws.Range("D9")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "25.73"
$c.Style = $s
$ws.Range("E9").Value = "  -6.50%  "
$c = $ws.Range("D10")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.06899"
$c.Style = $s
$ws.Range("E10").Value = "  -3.06%  "
$c = $ws.Range("D11")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.08018"
$c.Style = $s
$ws.Range("E11").Value = "  -0.41%  "
$c = $ws.Range("D12")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.7569"
$c.Style = $s
$ws.Range("E12").Value = "  -3.06%  "
$ws.Range("D13").Value = "1.906.36"
$ws.Range("E13").Value = "  -0.72%  "
$c = $ws.Range("D14")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "5.235"
$c.Style = $s
$ws.Range("E14").Value = "  -3.16%  "
$c = $ws.Range("D15")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "91.15"
$c.Style = $s
$ws.Range("E15").Value = "  -2.29%  "
$c = $ws.Range("D16")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "6.147"
$c.Style = $s
$ws.Range("E16").Value = "  +1.62%  "
$ws.Range("D17").Value = "30.033.12"
$ws.Range("E17").Value = "  -0.94%  "
$c = $ws.Range("D18")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "14.05"
$c.Style = $s
$ws.Range("E18").Value = "  -3.74%  "
$c = $ws.Range("D19")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.000007767"
$c.Style = $s
$ws.Range("E19").Value = "  -2.09%  "
$c = $ws.Range("D20")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "236.62"
$c.Style = $s
$ws.Range("E20").Value = "  -6.21%  "
$ws.Range("D21").Value = "2.160.37"
$ws.Range("E21").Value = "  -2.37%  "
$c = $ws.Range("D22")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.9995"
$c.Style = $s
$ws.Range("E22").Value = "  -0.04%  "
$ws.Range("E23").Value = "  -0.11%  "
$c = $ws.Range("D24")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "7.074"
$c.Style = $s
$ws.Range("E24").Value = "  +6.19%  "
$ws.Range("B25").Value = "Monero"
$ws.Range("C25").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$c = $ws.Range("D25")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "167.39"
$c.Style = $s
$ws.Range("E25").Value = "  +1.24%  "
$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$c = $ws.Range("D26")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "9.311"
$c.Style = $s
$ws.Range("E26").Value = "  -2.90%  "
$c = $ws.Range("D27")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "18.82"
$c.Style = $s
$ws.Range("E27").Value = "  -1.38%  "
$c = $ws.Range("D28")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.1265"
$c.Style = $s
$ws.Range("E28").Value = "  -2.14%  "
$c = $ws.Range("D29")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "2.054"
$c.Style = $s
$ws.Range("E29").Value = "  -6.12%  "
$c = $ws.Range("D30")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "1.351"
$c.Style = $s
$ws.Range("E30").Value = "  -0.78%  "
$ws.Range("E31").Value = "  -2.42%  "
$c = $ws.Range("D32")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "4.294"
$c.Style = $s
$ws.Range("E32").Value = "  -3.02%  "
$c = $ws.Range("D33")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "4.038"
$c.Style = $s
$ws.Range("E33").Value = "  -2.50%  "
$c = $ws.Range("D34")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.05324"
$c.Style = $s
$ws.Range("E34").Value = "  +1.37%  "
$c = $ws.Range("D35")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "1.287"
$c.Style = $s
$ws.Range("E35").Value = "  -2.47%  "
$c = $ws.Range("D36")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.7404"
$c.Style = $s
$c = $ws.Range("D37")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "2.722"
$c.Style = $s
$ws.Range("E37").Value = "  -1.59%  "
$c = $ws.Range("D38")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.01944"
$c.Style = $s
$ws.Range("E38").Value = "  -0.53%  "
$c = $ws.Range("D39")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "2.764"
$c.Style = $s
$ws.Range("E39").Value = "  -1.26%  "
$c = $ws.Range("D40")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "6.248"
$c.Style = $s
$ws.Range("E40").Value = "  -4.19%  "
$c = $ws.Range("D41")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.4463"
$c.Style = $s
$ws.Range("E41").Value = "  -1.45%  "
$c = $ws.Range("D42")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "72.66"
$c.Style = $s
$ws.Range("E42").Value = "  -6.12%  "
$c = $ws.Range("D43")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "1.956"
$c.Style = $s
$ws.Range("E43").Value = "  -0.77%  "
$c = $ws.Range("D44")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.9998"
$c.Style = $s
$ws.Range("E44").Value = "  -0.05%  "
$c = $ws.Range("D45")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.8310"
$c.Style = $s
$ws.Range("E45").Value = "  -1.60%  "
$c = $ws.Range("D46")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "7.705"
$c.Style = $s
$ws.Range("E46").Value = "  -0.12%  "
$c = $ws.Range("D47")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "101.32"
$c.Style = $s
$ws.Range("E47").Value = "  -0.13%  "
$c = $ws.Range("D48")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "9.826"
$c.Style = $s
$ws.Range("E48").Value = "  -2.22%  "
$ws.Range("D49").Value = "2.061.14"
$ws.Range("E49").Value = "  -1.60%  "
$c = $ws.Range("D50")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "36.60"
$c.Style = $s
$ws.Range("E50").Value = "  -3.32%  "
$c = $ws.Range("D51")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.1166"
$c.Style = $s
$ws.Range("E51").Value = "  -4.60%  "
